# reviewerTodoList.xlsx update
# Marks several additional TODO rows as COMPLETE (D column = TRUE) and
# records reviewer notes in the "notes" column (E), matching the author's
# responses to the thesis reviewer's comments. Overall completion moves
# from 40.74% to 62.96%.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reviewer-note text is written first, in the order the notes were
# originally authored, so new shared-string entries line up with the
# canonical workbook.
$ws.Cells.Item(15, 5).Value = "is this neccesary? I cannot remember every search word I used"
$ws.Cells.Item(47, 5).Value = "is this required?"
$ws.Cells.Item(54, 5).Value = "references include page numbers, is this required?"
$ws.Cells.Item(5, 5).Value = "I cannot find any parts where I do this"
$ws.Cells.Item(9, 5).Value = "not having yes answers makes the contributions less than they are"
$ws.Cells.Item(41, 5).Value = "Yes, but it was meant to point there, but I can understand the confusion"

# Row 53 reuses the same note as row 47 ("is this required?").
$ws.Cells.Item(53, 5).Value = "is this required?"
# Row 42 reuses the same note as row 9.
$ws.Cells.Item(42, 5).Value = "not having yes answers makes the contributions less than they are"

# Mark the corresponding TODO rows as COMPLETE.
# Row 5: "Abstract: note the rationale behind the qualitative and quantitative experiments"
$ws.Cells.Item(5, 4).Value = $true

# Row 6: "Intro: ground the thesis in a clear statement of contribution..."
$ws.Cells.Item(6, 4).Value = $true

# Row 7: "Intro: Number the research questions and separate them from the text"
$ws.Cells.Item(7, 4).Value = $true

# Row 9: "Intro: separate "research aims and contributions" into ..."
$ws.Cells.Item(9, 4).Value = $true

# Row 15: "LitRev: explain which databases were searched and note search words used"
$ws.Cells.Item(15, 4).Value = $true

# Row 20: "LitRev: ground this section in terms of the research questions more"
$ws.Cells.Item(20, 4).Value = $true

# Row 41: "Conclusions: Link results back to research questions"
$ws.Cells.Item(41, 4).Value = $true

# Row 42: "Conclusions: a cross reference to the literature review section..."
$ws.Cells.Item(42, 4).Value = $true

# Row 43: "Conclusions: answer research questions better than "yes" "no""
$ws.Cells.Item(43, 4).Value = $true

# Row 47: "Section 3.1: Explain better for non-technical readers"
$ws.Cells.Item(47, 4).Value = $true

# Row 53: "references should be of format: page numbers, dashed page numbers..."
$ws.Cells.Item(53, 4).Value = $true

# Row 54: "citations should include page numbers"
$ws.Cells.Item(54, 4).Value = $true

# Rows that now carry a two-line wrapped note grow to fit the text.
$ws.Rows.Item(9).RowHeight = 28.8
$ws.Rows.Item(15).RowHeight = 28.8
$ws.Rows.Item(41).RowHeight = 28.8
$ws.Rows.Item(42).RowHeight = 28.8

# Move the view/selection to reflect where the author was last working.
$null = $ws.Range("D34").Select()
